$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 468 (shifts existing rows 468-572 down to 469-573)
$ws.Rows('468:468').Insert()

# Fill in the new row's data
$ws.Range('A468').Value = 10
$ws.Range('B468').Value = 'Vega Modelo de Temuco'
$ws.Range('C468').Value = 'La Araucanía'
$ws.Range('D468').Value = 44642
$ws.Range('E468').Value = 9
$ws.Range('F468').Value = 'Fruta'
$ws.Range('G468').Value = 100103
$ws.Range('H468').Value = 'Frutos de hueso (carozo)'
$ws.Range('I468').Value = 100103006
$ws.Range('J468').Value = 'Nectarín'
$ws.Range('K468').Value = 'Artic Snow'
$ws.Range('L468').Value = 'Especial'
$ws.Range('M468').Value = 200
$ws.Range('N468').Value = 19000
$ws.Range('O468').Value = 19000
$ws.Range('P468').Value = 19000
$ws.Range('Q468').Value = '$/caja 15 kilos empedrada'
$ws.Range('R468').Value = "Región de O'Higgins"
$ws.Range('S468').Value = 1267
$ws.Range('T468').Value = 15
